$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "MaxInfected" section mirrors the existing Diversity/Peak.Strains/Total.Strains blocks:
# a title row, a header row (R0 = ... columns), and one data row.
$ws.Range("A25").Value = "MaxInfected"

$ws.Range("A26").Value = "R0 = 1.4"
$ws.Range("B26").Value = "R0 = 2.1"
$ws.Range("C26").Value = "R0 = 2.8"
$ws.Range("D26").Value = "R0 = 3.5"
$ws.Range("E26").Value = "R0 = 4.2"
$ws.Range("F26").Value = "R0 = 4.9"
$ws.Range("G26").Value = "R0 = 5.6"
$ws.Range("H26").Value = "R0 = 6.3"
$ws.Range("I26").Value = "R0 = 7"

$ws.Range("A27").Value = 0.044180610000000002
$ws.Range("B27").Value = 0.15590329999999999
$ws.Range("C27").Value = 0.27389229999999998
$ws.Range("D27").Value = 0.35549989999999998
$ws.Range("E27").Value = 0.41962870000000002
$ws.Range("F27").Value = 0.46964280000000003
$ws.Range("G27").Value = 0.51334049999999998
$ws.Range("H27").Value = 0.54596710000000004
$ws.Range("I27").Value = 0.57680920000000002

$ws.Range("A51").Select()
